$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 265.375
$ws.Range("I4").Value = 265.375
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 265.375
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -151.375
$ws.Range("N4").Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 294.16666
$ws.Range("I9").Value = 318.4
$ws.Range("K9").Value = 318.4
$ws.Range("M9").Value = -149.4
$ws.Range("N9").Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 247.11765
$ws.Range("I33").Value = 247.11765
$ws.Range("K33").Value = 247.11765
$ws.Range("M33").Value = -18.11765

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1179.3077
$ws.Range("I53").Value = 1212.3636
$ws.Range("J53").Value = 997.5
$ws.Range("K53").Value = 1212.3636
$ws.Range("L53").Value = 997.5
$ws.Range("M53").Value = -575.3635999999999
$ws.Range("N53").Value = -2271.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6333.3335
$ws.Range("I113").Value = 7000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 7000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -3746
$ws.Range("N113").Value = -11508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2252.25
$ws.Range("I138").Value = 1504.5
$ws.Range("K138").Value = 4513.5
$ws.Range("M138").Value = 626.5
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 867.9167
$ws.Range("I61").Value = 765
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 765
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -553
$ws.Range("N61").Value = -2424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 867.9167
$ws.Range("I136").Value = 765
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2295
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 255
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 20815.2
$ws.Range("I105").Value = 976.2857
$ws.Range("J105").Value = 67106
$ws.Range("K105").Value = 976.2857
$ws.Range("L105").Value = 67106
$ws.Range("M105").Value = 770.7143
$ws.Range("N105").Value = -70600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2166.923
$ws.Range("I31").Value = 735.25
$ws.Range("K31").Value = 735.25
$ws.Range("M31").Value = -440.25
$ws.Range("N31").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2166.923
$ws.Range("I34").Value = 735.25
$ws.Range("K34").Value = 735.25
$ws.Range("M34").Value = -533.25
$ws.Range("N34").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2318.375
$ws.Range("I58").Value = 1006.8
$ws.Range("K58").Value = 1006.8
$ws.Range("M58").Value = -803.8
$ws.Range("N58").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6274.1816
$ws.Range("J99").Value = 6303.7144
$ws.Range("L99").Value = 6303.7144
$ws.Range("N99").Value = -9299.714400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 604
$ws.Range("I122").Value = 618.8570999999999
$ws.Range("K122").Value = 1856.5713
$ws.Range("M122").Value = 593.4287000000002
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6274.1816
$ws.Range("J126").Value = 6303.7144
$ws.Range("L126").Value = 18911.1432
$ws.Range("N126").Value = -23851.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2318.375
$ws.Range("I136").Value = 1006.8
$ws.Range("K136").Value = 3020.4
$ws.Range("M136").Value = -470.3999999999996
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 183.125
$ws.Range("J12").Value = 285.66666
$ws.Range("L12").Value = 856.9999799999999
$ws.Range("N12").Value = -1202.99998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12099.875
$ws.Range("I80").Value = 2379.8
$ws.Range("K80").Value = 2379.8
$ws.Range("M80").Value = -1381.8
$ws.Range("N80").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 12099.875
$ws.Range("I83").Value = 2379.8
$ws.Range("K83").Value = 11899
$ws.Range("M83").Value = -6907
$ws.Range("N83").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3180
$ws.Range("I122").Value = 3225
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9675
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -7225
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6496.5
$ws.Range("I7").Value = 6494
$ws.Range("K7").Value = 6494
$ws.Range("M7").Value = -6382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 5000
$ws.Range("I26").Value = 5000
$ws.Range("K26").Value = 5000
$ws.Range("M26").Value = -4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1947
$ws.Range("I82").Value = 2122.3333
$ws.Range("J82").Value = 1771.6666
$ws.Range("K82").Value = 2122.3333
$ws.Range("L82").Value = 1771.6666
$ws.Range("M82").Value = -1761.3333
$ws.Range("N82").Value = -2493.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1947
$ws.Range("I85").Value = 2122.3333
$ws.Range("J85").Value = 1771.6666
$ws.Range("K85").Value = 2122.3333
$ws.Range("L85").Value = 1771.6666
$ws.Range("M85").Value = -874.3332999999998
$ws.Range("N85").Value = -4267.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 993.75
$ws.Range("I93").Value = 991.6667
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 991.6667
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 256.3333
$ws.Range("N93").Value = -3496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 49000
$ws.Range("J104").Value = 49000
$ws.Range("L104").Value = 49000
$ws.Range("N104").Value = -55988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6496.5
$ws.Range("I126").Value = 6494
$ws.Range("K126").Value = 19482
$ws.Range("M126").Value = -17012
$ws.Range("N126").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 45193.5
$ws.Range("J112").Value = 45193.5
$ws.Range("L112").Value = 45193.5
$ws.Range("N112").Value = -48147.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1634.6666
$ws.Range("I122").Value = 1634.6666
$ws.Range("K122").Value = 4903.9998
$ws.Range("M122").Value = -2453.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3234.8
$ws.Range("I126").Value = 3234.8
$ws.Range("K126").Value = 9704.400000000001
$ws.Range("M126").Value = -7234.400000000001
